# Add team record (Wins/Losses/Ties) columns to the player data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF, matching the style of the
# existing header cells (bold, centered, thin border) via format copy.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record for every player row (2-52) with the same
# constant values (92 wins, 70 losses, 0 ties).
$lastRow = 52
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 92
    $ws.Cells.Item($r, 31).Value = 70
    $ws.Cells.Item($r, 32).Value = 0
}
